# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (variety "Dina") right before the existing
# row 24, pushing the former rows 24-33 down to rows 27-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at position 24 (shifts old rows 24-33 -> 27-36,
# and grows the sheet dimension from A1:T33 to A1:T36 automatically).
$ws.Rows("24:26").Insert()

# Columns shared by every record in this data block.
$ws.Range("A24:A26").Value = 2
$ws.Range("B24:B26").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C24:C26").Value = "Coquimbo"
$ws.Range("D24:D26").Value = 44924
$ws.Range("E24:E26").Value = 4
$ws.Range("F24:F26").Value = "Fruta"
$ws.Range("G24:G26").Value = 100103
$ws.Range("H24:H26").Value = "Frutos de hueso (carozo)"
$ws.Range("I24:I26").Value = 100103003
$ws.Range("J24:J26").Value = "Damasco"
$ws.Range("K24:K26").Value = "Dina"
$ws.Range("Q24:Q26").Value = "$/caja 16 kilos"
$ws.Range("R24:R26").Value = "Región Metropolitana"

# Row 24: Dina / Especial
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 200
$ws.Range("N24").Value = 23000
$ws.Range("O24").Value = 24000
$ws.Range("P24").Value = 23500
$ws.Range("S24").Value = 1469
$ws.Range("T24").Value = 16

# Row 25: Dina / Primera
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 300
$ws.Range("N25").Value = 20000
$ws.Range("O25").Value = 21000
$ws.Range("P25").Value = 20500
$ws.Range("S25").Value = 1281
$ws.Range("T25").Value = 16

# Row 26: Dina / Segunda
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 300
$ws.Range("N26").Value = 15000
$ws.Range("O26").Value = 16000
$ws.Range("P26").Value = 15500
$ws.Range("S26").Value = 969
$ws.Range("T26").Value = 16
